$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, column letter, new text value.
# Columns D (Price) and E (Volume 1h) are plain display strings (with
# thousands dots / leading+trailing spaces / percent signs) that must be
# stored verbatim as text, so those are written through a short
# "force text" dance: set NumberFormat to "@" (Text) before assigning the
# value so Excel does not reinterpret e.g. "1.00" or "18.03" as a number,
# then restore the cell style to "Normal" so no stray formatting remains.
$updates = @(
    @{Row=2; Col="D"; Value="69.606.85"}
    @{Row=2; Col="E"; Value="  -0.62%  "}
    @{Row=3; Col="D"; Value="3.546.40"}
    @{Row=3; Col="E"; Value="  -1.97%  "}
    @{Row=4; Col="D"; Value="1.00"}
    @{Row=4; Col="E"; Value="  +0.31%  "}
    @{Row=5; Col="D"; Value="198.17"}
    @{Row=5; Col="E"; Value="  +1.42%  "}
    @{Row=6; Col="D"; Value="586.17"}
    @{Row=6; Col="E"; Value="  -3.20%  "}
    @{Row=7; Col="D"; Value="0.614"}
    @{Row=7; Col="E"; Value="  -1.98%  "}
    @{Row=8; Col="E"; Value="  +0.06%  "}
    @{Row=9; Col="E"; Value="  +0.24%  "}
    @{Row=10; Col="D"; Value="0.630"}
    @{Row=10; Col="E"; Value="  -3.20%  "}
    @{Row=11; Col="D"; Value="52.16"}
    @{Row=11; Col="E"; Value="  -3.36%  "}
    @{Row=12; Col="E"; Value="  -4.88%  "}
    @{Row=13; Col="B"; Value="Polkadot"}
    @{Row=13; Col="C"; Value="https://coinranking.com/coin/25W7FG7om+polkadot-dot"}
    @{Row=13; Col="D"; Value="9.34"}
    @{Row=13; Col="E"; Value="  -2.18%  "}
    @{Row=14; Col="B"; Value="BitcoinCash"}
    @{Row=14; Col="C"; Value="https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"}
    @{Row=14; Col="D"; Value="687.47"}
    @{Row=14; Col="E"; Value="  +16.15%  "}
    @{Row=15; Col="D"; Value="4.109.09"}
    @{Row=15; Col="E"; Value="  -1.91%  "}
    @{Row=16; Col="D"; Value="69.684.42"}
    @{Row=16; Col="E"; Value="  -0.79%  "}
    @{Row=17; Col="D"; Value="3.555.24"}
    @{Row=17; Col="E"; Value="  -1.88%  "}
    @{Row=18; Col="D"; Value="12.49"}
    @{Row=18; Col="E"; Value="  -5.94%  "}
    @{Row=19; Col="D"; Value="18.61"}
    @{Row=19; Col="E"; Value="  -3.12%  "}
    @{Row=20; Col="E"; Value="  -0.73%  "}
    @{Row=21; Col="D"; Value="0.972"}
    @{Row=21; Col="E"; Value="  -2.21%  "}
    @{Row=22; Col="D"; Value="18.03"}
    @{Row=23; Col="D"; Value="108.44"}
    @{Row=23; Col="E"; Value="  +5.69%  "}
    @{Row=24; Col="D"; Value="5.22"}
    @{Row=24; Col="E"; Value="  +0.80%  "}
    @{Row=25; Col="D"; Value="4.42"}
    @{Row=25; Col="E"; Value="  -4.54%  "}
    @{Row=26; Col="E"; Value="  -2.92%  "}
    @{Row=27; Col="E"; Value="  -4.22%  "}
    @{Row=28; Col="D"; Value="9.72"}
    @{Row=28; Col="E"; Value="  +1.09%  "}
    @{Row=29; Col="D"; Value="33.69"}
    @{Row=29; Col="E"; Value="  -0.60%  "}
    @{Row=30; Col="E"; Value="  -1.34%  "}
    @{Row=31; Col="D"; Value="6.91"}
    @{Row=31; Col="E"; Value="  -2.89%  "}
    @{Row=32; Col="D"; Value="11.95"}
    @{Row=32; Col="E"; Value="  -3.22%  "}
    @{Row=33; Col="E"; Value="  -4.07%  "}
    @{Row=34; Col="D"; Value="62.05"}
    @{Row=34; Col="E"; Value="  -1.86%  "}
    @{Row=35; Col="D"; Value="3.811.21"}
    @{Row=35; Col="E"; Value="  -3.58%  "}
    @{Row=36; Col="D"; Value="0.0₃0821"}
    @{Row=36; Col="E"; Value="  -9.24%  "}
    @{Row=37; Col="E"; Value="  +0.13%  "}
    @{Row=38; Col="E"; Value="  +4.40%  "}
    @{Row=39; Col="B"; Value="Fetch.AI"}
    @{Row=39; Col="C"; Value="https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"}
    @{Row=39; Col="D"; Value="2.95"}
    @{Row=39; Col="E"; Value="  -6.52%  "}
    @{Row=40; Col="B"; Value="Bittensor"}
    @{Row=40; Col="C"; Value="https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"}
    @{Row=40; Col="D"; Value="499.01"}
    @{Row=40; Col="E"; Value="  -5.63%  "}
    @{Row=41; Col="E"; Value="  -4.63%  "}
    @{Row=42; Col="E"; Value="  +1.53%  "}
    @{Row=43; Col="D"; Value="34.76"}
    @{Row=43; Col="E"; Value="  -6.57%  "}
    @{Row=44; Col="E"; Value="  +0.96%  "}
    @{Row=45; Col="D"; Value="2.94"}
    @{Row=45; Col="E"; Value="  +2.86%  "}
    @{Row=46; Col="E"; Value="  +1.22%  "}
    @{Row=47; Col="E"; Value="  -2.34%  "}
    @{Row=48; Col="E"; Value="  -0.34%  "}
    @{Row=49; Col="D"; Value="8.43"}
    @{Row=49; Col="E"; Value="  -2.26%  "}
    @{Row=50; Col="B"; Value="Jupiter"}
    @{Row=50; Col="C"; Value="https://coinranking.com/coin/qMgTxtv34+jupiter-jup"}
    @{Row=50; Col="D"; Value="1.82"}
    @{Row=50; Col="E"; Value="  +22.17%  "}
    @{Row=51; Col="B"; Value="CoreDAO"}
    @{Row=51; Col="C"; Value="https://coinranking.com/coin/HFvoXUQh4+coredao-core"}
    @{Row=51; Col="D"; Value="2.81"}
    @{Row=51; Col="E"; Value="  +70.30%  "}
)

foreach ($u in $updates) {
    $cell = $ws.Range("$($u.Col)$($u.Row)")
    if ($u.Col -eq "D" -or $u.Col -eq "E") {
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}
